$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated ticker grid for rows 2-17 (columns B:F). Values of $null clear the cell.
$data = @{
    2  = @("NSE:ADANIENSOL", "NSE:ASKAUTOLTD", "NSE:AMBUJACEM", "NSE:ASTRAL",  "NSE:ADANIENSOL")
    3  = @("NSE:AETHER",     "NSE:BAJAJHCARE", $null,           "NSE:BSE",     "NSE:BLUESTARCO")
    4  = @("NSE:BLUESTARCO", "NSE:BASML",      $null,           "NSE:DLF",     "NSE:CHOLAFIN")
    5  = @("NSE:CHOLAFIN",   "NSE:BFUTILITIE", $null,           "NSE:EXIDEIND","NSE:HEROMOTOCO")
    6  = @("NSE:DIL",        "NSE:DCMSHRIRAM", $null,           $null,         "NSE:LICI")
    7  = @("NSE:DREDGECORP", "NSE:EMAMIPAP",   $null,           $null,         $null)
    8  = @("NSE:GHCL",       "NSE:IPL",        $null,           $null,         $null)
    9  = @("NSE:HEIDELBERG", "NSE:JUBLFOOD",   $null,           $null,         $null)
    10 = @("NSE:HEROMOTOCO", "NSE:KCP",        $null,           $null,         $null)
    11 = @("NSE:LICI",       "NSE:KIRIINDUS",  $null,           $null,         $null)
    12 = @("NSE:MAHAPEXLTD", "NSE:LUXIND",     $null,           $null,         $null)
    13 = @("NSE:MMTC",       "NSE:MSUMI",      $null,           $null,         $null)
    14 = @("NSE:NIFMID150",  "NSE:MVGJL",      $null,           $null,         $null)
    15 = @($null,            "NSE:NEULANDLAB", $null,           $null,         $null)
    16 = @($null,            "NSE:PIDILITIND", $null,           $null,         $null)
    17 = @($null,            "NSE:REPL",       $null,           $null,         $null)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 2 + $i   # B=2 .. F=6
        $cell = $ws.Cells.Item($r, $col)
        if ($vals[$i] -eq $null) {
            $cell.Value = ""
        } else {
            $cell.Value = $vals[$i]
        }
    }
}

# Rows 18 and 19 (NSE:ORCHPHARMA, NSE:PNBHOUSING) are dropped entirely.
$ws.Rows(18).Delete()
$ws.Rows(18).Delete()
